$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for rows 2-51 to reflect refreshed crypto data.
# Numeric-looking Price values are forced to remain plain text (matching original inlineStr cells)
# by briefly applying a text NumberFormat, then restoring the default "Normal" style.

$ws.Range("D2").Value = '41.363.46'
$ws.Range("E2").Value = '  +0.76%  '
$ws.Range("D3").Value = '2.190.03'
$ws.Range("E3").Value = '  +0.06%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '254.56'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.72%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.617'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.39%  '
$ws.Range("E7").Value = '  -1.94%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.586'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +9.16%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.11'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '58.73'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.81%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0940'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.41%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.20'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +10.06%  '
$ws.Range("E14").Value = '  +0.74%  '
$ws.Range("D15").Value = '2.515.65'
$ws.Range("E15").Value = '  +0.16%  '
$ws.Range("E16").Value = '  +5.21%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.53'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.07%  '
$ws.Range("D18").Value = '2.184.15'
$ws.Range("E18").Value = '  -0.18%  '
$ws.Range("D19").Value = '41.292.48'
$ws.Range("E19").Value = '  +0.74%  '
$ws.Range("D20").Value = '0.0₃0956'
$ws.Range("E20").Value = '  +2.04%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.23'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.59%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.18'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.14%  '
$ws.Range("E23").Value = '  +1.18%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.05'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.93%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.91'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +23.22%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.91'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +8.87%  '
$ws.Range("E27").Value = '  +0.02%  '
$ws.Range("E28").Value = '  +5.34%  '
$ws.Range("E29").Value = '  -0.38%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '170.20'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.29%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.70'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.86%  '
$ws.Range("E32").Value = '  +2.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.55'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +9.56%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.123'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.67%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0740'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.22%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '26.78'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +14.91%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.63'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.13'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +7.96%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0302'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +13.87%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '12.83'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +27.87%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.21'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.90%  '
$ws.Range("E42").Value = '  -1.71%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '64.63'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.19%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.02'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.90%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.203'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.74%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.64'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.57%  '
$ws.Range("E47").Value = '  +3.90%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.01'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.24%  '
$ws.Range("E49").Value = '  +5.52%  '
$ws.Range("E50").Value = '  +1.63%  '
$ws.Range("E51").Value = '  -3.80%  '
